# arreglo error formato excel
$wb = $excel.ActiveWorkbook

# --- Sheet "Administrativos" ---
$ws1 = $wb.Worksheets.Item("Administrativos")
$ws1.Range("A2").Value = "Nombres 1"
$ws1.Range("C2").Value = "identificacion1"

# --- Sheet "Docentes-Conciliadores" ---
$ws2 = $wb.Worksheets.Item("Docentes-Conciliadores")
$ws2.Range("A2").Value = "Nombres 2"
$ws2.Range("C2").Value = "identificacion2"
$ws2.Range("A3").Value = "Nombres 3"

# --- Sheet "Estudiantes" ---
$ws3 = $wb.Worksheets.Item("Estudiantes")
$ws3.Range("A2").Value = "Nombres 4"
$ws3.Range("C2").Value = 202215203458
$ws3.Range("A3").Value = "Nombres 5"
$ws3.Range("A4").Value = "Nombres 6"
$ws3.Range("A5").Value = "Nombres 7"
$ws3.Range("A6").Value = "Nombres 8"
$ws3.Range("A7").Value = "Nombres 9"
$ws3.Range("A8").Value = "Nombres 10"
$ws3.Range("A9").Value = "Nombres 11"
$ws3.Range("A10").Value = "Nombres 12"

# --- back to Docentes-Conciliadores for the identificacion3 value ---
$ws2.Range("C3").Value = "identificacion3"
